# Generate Report for Handoff
# Swap the handed-off source file's identifier (UUID-based file name) and
# refresh the associated handoff/handback timestamps across the three
# sheets (Overview, zh-cn, de-de) of the localization-status workbook.

$wb = $excel.ActiveWorkbook

$oldGuid = "006b9da5-e77f-45e3-889c-00792bc230f2"
$newGuid = "731d9710-f734-40b9-a4c7-1ba3d02abd22"

$oldMd = "$oldGuid.md"
$newMd = "$newGuid.md"

$oldMdPath = "e2e\$oldMd"
$newMdPath = "e2e\$newMd"

$oldZhXlf = "$oldGuid.5ba8fcea68b1a15e0faf1051e18426a77ac24608.zh-cn.xlf"
$newZhXlf = "$newGuid.a94eec8ea9c635f9b538420a3a322b7e4a73803c.zh-cn.xlf"

$oldDeXlf = "$oldGuid.5ba8fcea68b1a15e0faf1051e18426a77ac24608.de-de.xlf"
$newDeXlf = "$newGuid.a94eec8ea9c635f9b538420a3a322b7e4a73803c.de-de.xlf"

function Update-HyperlinkDisplay {
    param($range, [string]$address, [string]$displayText)

    # The hyperlink's underlying target address is untouched by this edit -
    # only the human-readable display text changes. Recreate the hyperlink
    # in place (delete + re-add) so the <hyperlink> element keeps pointing
    # at the same address while showing the new text.
    # (Note: use positional args - named "-param value" binding isn't
    # reliable for custom functions in this host.)
    $range.Hyperlinks.Delete()
    $range.Worksheet.Hyperlinks.Add($range, $address, "", "", $displayText) | Out-Null
}

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$overviewAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/351da6a858b3aafa85f045f833a5438c152227ff/e2e/$oldMd"

$wsOverview.Range("A2").Value = $newMd
Update-HyperlinkDisplay $wsOverview.Range("B2") $overviewAddress $newMdPath
$wsOverview.Range("G2").Value = "2016-08-30 09:04:40"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$zhCnAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/351da6a858b3aafa85f045f833a5438c152227ff/e2e/$oldMd"

Update-HyperlinkDisplay $wsZhCn.Range("A2") $zhCnAddress $newMd
$wsZhCn.Range("G2").Value = $newZhXlf
$wsZhCn.Range("H2").Value = "2016-08-30 09:04:29"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$deDeAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/351da6a858b3aafa85f045f833a5438c152227ff/e2e/$oldMd"

Update-HyperlinkDisplay $wsDeDe.Range("A2") $deDeAddress $newMd
$wsDeDe.Range("G2").Value = $newDeXlf
$wsDeDe.Range("H2").Value = "2016-08-30 09:04:40"
